# Update the "modtage" (receive) overview sheet for the new publish date.
# The workbook tracks its own "last updated" date in the sheet's name, and
# the defined name `Privathospitalssystemer___modtage` simply points at that
# sheet, so renaming the sheet is the one edit required here - Excel keeps
# the defined name's formula in sync with the sheet name automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Name = "Opdateret d. 05-12-2025"
